$d = $word.ActiveDocument

# The title ("Računalna forenzika: Laboratorijske vježbe 1") lives in the
# very first run of the very first paragraph.
$titlePara  = $d.Paragraphs(1).Range
$titleStart = $titlePara.Start
$oldTitle   = "Računalna forenzika: Laboratorijske vježbe 1"
$splitLen   = "Računalna fo".Length
$titleEnd   = $titleStart + $oldTitle.Length

# 1) Bump the lab number in the title: "vježbe 1" -> "vježbe 3".
$d.Content.Find.Execute("Laboratorijske vježbe 1", $true, $false, $false, $false, $false, $true, 1, $false, "Laboratorijske vježbe 3", 2)

# 2) Break the (now updated) title into two runs - "Računalna fo" and
#    "renzika: Laboratorijske vježbe 3" - by round-tripping a character
#    formatting property on just the first chunk. Changing and then
#    restoring Font.Size forces Word to split the run at that boundary
#    without altering any visible formatting.
$splitRange = $d.Range($titleStart, $titleStart + $splitLen)
$splitRange.Font.Size = $splitRange.Font.Size + 1
$splitRange.Font.Size = $splitRange.Font.Size - 1

# 3) Re-plant the "_GoBack" bookmark right at the new run boundary, i.e.
#    between the two title runs and before the following tab character.
#    A document can only have one bookmark of a given name, so adding it
#    here automatically removes/relocates the bookmark that used to sit
#    in front of the picture in the second paragraph.
$bkPoint = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $bkPoint)
